# Updates cryptos list values (price/volume columns) per the data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a text value even when it looks numeric,
    # matching the inline-string / text representation used in the source file,
    # then restore the default "Normal" style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "38.796.26"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "2.104.87"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "227.78"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +0.55%  "
Set-TextValue $ws.Range("D7") "62.60"
$ws.Range("E7").Value = "  +2.88%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +2.29%  "
Set-TextValue $ws.Range("D10") "0.0844"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  +6.42%  "
$ws.Range("D13").Value = "2.416.62"
$ws.Range("E13").Value = "  +0.67%  "
Set-TextValue $ws.Range("D14") "22.04"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "2.102.80"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "38.804.21"
$ws.Range("E18").Value = "  +1.23%  "
Set-TextValue $ws.Range("D19") "6.12"
$ws.Range("E19").Value = "  +1.01%  "
Set-TextValue $ws.Range("D20") "71.53"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  +1.72%  "
Set-TextValue $ws.Range("D22") "228.24"
$ws.Range("E22").Value = "  +1.21%  "
Set-TextValue $ws.Range("D24") "2.36"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("E25").Value = "  +0.32%  "
Set-TextValue $ws.Range("D26") "9.66"
$ws.Range("E26").Value = "  +2.29%  "
Set-TextValue $ws.Range("D27") "172.46"
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("E28").Value = "  +1.13%  "
Set-TextValue $ws.Range("D29") "1.41"
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("E30").Value = "  +1.63%  "
Set-TextValue $ws.Range("D31") "2.58"
$ws.Range("E31").Value = "  +10.65%  "
$ws.Range("E32").Value = "  +0.27%  "
Set-TextValue $ws.Range("D33") "4.58"
$ws.Range("E33").Value = "  +1.55%  "
Set-TextValue $ws.Range("D34") "7.18"
$ws.Range("E34").Value = "  +11.80%  "
Set-TextValue $ws.Range("D35") "4.76"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  -0.42%  "
Set-TextValue $ws.Range("D39") "0.999"
$ws.Range("E39").Value = "  -0.22%  "
Set-TextValue $ws.Range("D40") "18.08"
$ws.Range("E40").Value = "  -1.98%  "
Set-TextValue $ws.Range("D41") "102.40"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("D43").Value = "1.526.68"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("E44").Value = "  +7.66%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D45") "2.81"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D46") "0.0917"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D47") "7.79"
$ws.Range("E47").Value = "  +0.68%  "
Set-TextValue $ws.Range("D48") "4.19"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "2.303.57"
$ws.Range("E51").Value = "  +0.72%  "
